$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 8: update D8, E8 (male group) and J8 (female group)
$ws.Range("D8").Value = 0.82
$ws.Range("E8").Value = 1.018
$ws.Range("J8").Value = 1.018

# Row 15: reset D15, E15 (male group) and J15 (female group) back to 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("J15").Value = 1

# Row 16: fill in missing J16 value
$ws.Range("J16").Value = 1

# Remove now-unused helper/placeholder cells in columns O:AA for rows 5-6
$ws.Range("O5:AA5").Clear()
$ws.Range("AA6").Clear()

# Update selection to match the saved view state
$ws.Range("J16").Select()
